$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 825
$ws1.Range("F6").Value = 137
$ws1.Range("F8").Value = 4884
$ws1.Range("F9").Value = 106
$ws1.Range("F10").Value = 5187
$ws1.Range("F11").Value = 592
$ws1.Range("F12").Value = 1299
$ws1.Range("F13").Value = 97

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 825
$ws4.Range("F6").Value = 137
$ws4.Range("F9").Value = 4884
$ws4.Range("F10").Value = 106
$ws4.Range("F11").Value = 5187
$ws4.Range("F12").Value = 592
$ws4.Range("F13").Value = 1299
$ws4.Range("F14").Value = 97
